$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1793103448275862
$ws.Range("C2").Value = 0.6137931034482759
$ws.Range("J2").Value = 0.02068965517241379
$ws.Range("P2").Value = 0.1275862068965517
$ws.Range("S2").Value = 0.05862068965517241

# Row 3
$ws.Range("B3").Value = 0.01075268817204301
$ws.Range("C3").Value = 0.03763440860215054
$ws.Range("J3").Value = 0.01075268817204301
$ws.Range("P3").Value = 0.7419354838709677
$ws.Range("S3").Value = 0.1989247311827957

# Row 4
$ws.Range("J4").Value = 0.02222222222222222
$ws.Range("P4").Value = 0.6888888888888889
$ws.Range("S4").Value = 0.2888888888888889

# Row 6
$ws.Range("B6").Value = 0.06103286384976526
$ws.Range("D6").Value = 0.009389671361502348
$ws.Range("F6").Value = 0.08450704225352113
$ws.Range("J6").Value = 0.2018779342723005
$ws.Range("O6").Value = 0.01408450704225352
$ws.Range("Q6").Value = 0.1924882629107981
$ws.Range("R6").Value = 0.09389671361502347
$ws.Range("S6").Value = 0.3427230046948357

# Row 7
$ws.Range("B7").Value = 0.1111111111111111
$ws.Range("D7").Value = 0.0308641975308642
$ws.Range("E7").Value = 0.006172839506172839
$ws.Range("F7").Value = 0.04938271604938271
$ws.Range("J7").Value = 0.09259259259259259
$ws.Range("Q7").Value = 0.1975308641975309
$ws.Range("R7").Value = 0.1296296296296296
$ws.Range("S7").Value = 0.382716049382716

# Row 8
$ws.Range("B8").Value = 0.08863636363636364
$ws.Range("D8").Value = 0.025
$ws.Range("E8").Value = 0.002272727272727273
$ws.Range("F8").Value = 0.06363636363636363
$ws.Range("J8").Value = 0.08636363636363636
$ws.Range("O8").Value = 0.02272727272727273
$ws.Range("Q8").Value = 0.1545454545454545
$ws.Range("R8").Value = 0.1204545454545455
$ws.Range("S8").Value = 0.4363636363636363

# Row 9
$ws.Range("B9").Value = 0.07843137254901961
$ws.Range("D9").Value = 0.01176470588235294
$ws.Range("E9").Value = 0.00392156862745098
$ws.Range("F9").Value = 0.0392156862745098
$ws.Range("J9").Value = 0.08235294117647059
$ws.Range("O9").Value = 0.0392156862745098
$ws.Range("Q9").Value = 0.1725490196078431
$ws.Range("R9").Value = 0.1058823529411765
$ws.Range("S9").Value = 0.4666666666666667

# Row 10
$ws.Range("B10").Value = 0.1176470588235294
$ws.Range("D10").Value = 0.01988400994200497
$ws.Range("E10").Value = 0.0008285004142502071
$ws.Range("F10").Value = 0.07456503728251865
$ws.Range("J10").Value = 0.1101905550952775
$ws.Range("O10").Value = 0.02236951118475559
$ws.Range("Q10").Value = 0.1971830985915493
$ws.Range("R10").Value = 0.09113504556752279
$ws.Range("S10").Value = 0.3661971830985916

# Row 11
$ws.Range("G11").Value = 0.1573426573426573
$ws.Range("J11").Value = 0.1363636363636364
$ws.Range("K11").Value = 0.2517482517482518
$ws.Range("L11").Value = 0.4265734265734266
$ws.Range("S11").Value = 0.02797202797202797

# Row 12
$ws.Range("G12").Value = 0.7479674796747967
$ws.Range("J12").Value = 0.1951219512195122
$ws.Range("L12").Value = 0.02439024390243903
$ws.Range("S12").Value = 0.03252032520325204

# Row 13
$ws.Range("G13").Value = 0.5849056603773585
$ws.Range("J13").Value = 0.3207547169811321
$ws.Range("S13").Value = 0.09433962264150944

# Row 14
$ws.Range("G14").Value = 0.25
$ws.Range("J14").Value = 0.25
$ws.Range("S14").Value = 0.5

# Row 15
$ws.Range("F15").Value = 0.008264462809917356
$ws.Range("H15").Value = 0.1611570247933884
$ws.Range("I15").Value = 0.08677685950413223
$ws.Range("J15").Value = 0.3223140495867768
$ws.Range("K15").Value = 0.04132231404958678
$ws.Range("M15").Value = 0.02066115702479339
$ws.Range("O15").Value = 0.09090909090909091
$ws.Range("S15").Value = 0.268595041322314

# Row 16
$ws.Range("F16").Value = 0.01990049751243781
$ws.Range("H16").Value = 0.1840796019900497
$ws.Range("I16").Value = 0.1144278606965174
$ws.Range("J16").Value = 0.3482587064676617
$ws.Range("K16").Value = 0.1194029850746269
$ws.Range("M16").Value = 0.02985074626865672
$ws.Range("N16").Value = 0.009950248756218905
$ws.Range("O16").Value = 0.05472636815920398
$ws.Range("S16").Value = 0.1194029850746269

# Row 17
$ws.Range("F17").Value = 0.00936768149882904
$ws.Range("H17").Value = 0.1733021077283372
$ws.Range("I17").Value = 0.1334894613583138
$ws.Range("J17").Value = 0.405152224824356
$ws.Range("K17").Value = 0.09601873536299765
$ws.Range("M17").Value = 0.01639344262295082
$ws.Range("O17").Value = 0.06088992974238876
$ws.Range("S17").Value = 0.1053864168618267

# Row 18
$ws.Range("F18").Value = 0.03478260869565217
$ws.Range("H18").Value = 0.1869565217391304
$ws.Range("I18").Value = 0.1043478260869565
$ws.Range("J18").Value = 0.391304347826087
$ws.Range("K18").Value = 0.08695652173913043
$ws.Range("M18").Value = 0.02173913043478261
$ws.Range("N18").Value = 0.008695652173913044
$ws.Range("O18").Value = 0.08695652173913043
$ws.Range("S18").Value = 0.0782608695652174

# Row 19
$ws.Range("F19").Value = 0.01037509976057462
$ws.Range("H19").Value = 0.1987230646448523
$ws.Range("I19").Value = 0.1037509976057462
$ws.Range("J19").Value = 0.3822825219473264
$ws.Range("K19").Value = 0.09177972865123703
$ws.Range("M19").Value = 0.02394253790901836
$ws.Range("O19").Value = 0.07023144453312051
$ws.Range("S19").Value = 0.1189146049481245
